$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Row=6; I='sd'; J='Statement-non-opinion'},
    @{Row=10; I='ba'; J='Appreciation'},
    @{Row=15; I='ba'; J='Appreciation'},
    @{Row=20; I='aa'; J='Agree/Accept'},
    @{Row=31; I='ba'; J='Appreciation'},
    @{Row=38; I='sd'; J='Statement-non-opinion'},
    @{Row=45; I='ba'; J='Appreciation'},
    @{Row=48; I='sd'; J='Statement-non-opinion'},
    @{Row=50; I='sd'; J='Statement-non-opinion'},
    @{Row=51; I='sv'; J='Statement-opinion'},
    @{Row=62; I='aa'; J='Agree/Accept'},
    @{Row=63; I='aa'; J='Agree/Accept'},
    @{Row=64; I='aa'; J='Agree/Accept'},
    @{Row=65; I='sd'; J='Statement-non-opinion'},
    @{Row=71; I='ba'; J='Appreciation'},
    @{Row=73; I='sv'; J='Statement-opinion'},
    @{Row=82; I='ba'; J='Appreciation'},
    @{Row=91; I='sd'; J='Statement-non-opinion'},
    @{Row=92; I='sd'; J='Statement-non-opinion'},
    @{Row=93; I='sd'; J='Statement-non-opinion'},
    @{Row=97; I='sd'; J='Statement-non-opinion'},
    @{Row=107; I='sd'; J='Statement-non-opinion'},
    @{Row=138; I='ba'; J='Appreciation'},
    @{Row=142; I='sd'; J='Statement-non-opinion'},
    @{Row=143; I='sd'; J='Statement-non-opinion'},
    @{Row=144; I='sd'; J='Statement-non-opinion'},
    @{Row=146; I='sd'; J='Statement-non-opinion'},
    @{Row=152; I='ba'; J='Appreciation'},
    @{Row=159; I='ba'; J='Appreciation'},
    @{Row=190; I='ba'; J='Appreciation'},
    @{Row=200; I='ba'; J='Appreciation'},
    @{Row=209; I='ba'; J='Appreciation'},
    @{Row=215; I='ba'; J='Appreciation'},
    @{Row=221; I='sd'; J='Statement-non-opinion'},
    @{Row=241; I='%'; J='Uninterpretable'},
    @{Row=250; I='ba'; J='Appreciation'},
    @{Row=252; I='ba'; J='Appreciation'},
    @{Row=255; I='ba'; J='Appreciation'},
    @{Row=264; I='aa'; J='Agree/Accept'},
    @{Row=267; I='%'; J='Uninterpretable'},
    @{Row=268; I='%'; J='Uninterpretable'},
    @{Row=278; I='ba'; J='Appreciation'},
    @{Row=281; I='sd'; J='Statement-non-opinion'},
    @{Row=284; I='aa'; J='Agree/Accept'},
    @{Row=286; I='ba'; J='Appreciation'},
    @{Row=287; I='%'; J='Uninterpretable'},
    @{Row=288; I='aa'; J='Agree/Accept'},
    @{Row=289; I='aa'; J='Agree/Accept'},
    @{Row=290; I='aa'; J='Agree/Accept'},
    @{Row=291; I='aa'; J='Agree/Accept'},
    @{Row=292; I='ba'; J='Appreciation'},
    @{Row=297; I='ba'; J='Appreciation'},
    @{Row=300; I='ba'; J='Appreciation'},
    @{Row=302; I='sv'; J='Statement-opinion'},
    @{Row=306; I='sd'; J='Statement-non-opinion'},
    @{Row=312; I='ba'; J='Appreciation'},
    @{Row=325; I='ba'; J='Appreciation'},
    @{Row=332; I='ba'; J='Appreciation'},
    @{Row=337; I='ba'; J='Appreciation'},
    @{Row=340; I='ba'; J='Appreciation'},
    @{Row=344; I='sd'; J='Statement-non-opinion'},
    @{Row=345; I='sd'; J='Statement-non-opinion'},
    @{Row=354; I='aa'; J='Agree/Accept'},
    @{Row=358; I='ba'; J='Appreciation'},
    @{Row=366; I='sd'; J='Statement-non-opinion'},
    @{Row=368; I='sv'; J='Statement-opinion'},
    @{Row=370; I='aa'; J='Agree/Accept'},
    @{Row=371; I='aa'; J='Agree/Accept'},
    @{Row=377; I='ba'; J='Appreciation'},
    @{Row=380; I='ba'; J='Appreciation'},
    @{Row=382; I='sd'; J='Statement-non-opinion'},
    @{Row=390; I='sv'; J='Statement-opinion'},
    @{Row=399; I='aa'; J='Agree/Accept'},
    @{Row=406; I='ba'; J='Appreciation'},
    @{Row=414; I='sv'; J='Statement-opinion'},
    @{Row=417; I='ba'; J='Appreciation'},
    @{Row=419; I='ba'; J='Appreciation'},
    @{Row=426; I='ba'; J='Appreciation'},
    @{Row=433; I='ba'; J='Appreciation'},
    @{Row=435; I='sd'; J='Statement-non-opinion'},
    @{Row=437; I='ba'; J='Appreciation'},
    @{Row=441; I='ba'; J='Appreciation'},
    @{Row=448; I='sv'; J='Statement-opinion'},
    @{Row=449; I='sd'; J='Statement-non-opinion'},
    @{Row=459; I='ba'; J='Appreciation'},
    @{Row=462; I='ba'; J='Appreciation'},
    @{Row=465; I='ba'; J='Appreciation'},
    @{Row=472; I='ba'; J='Appreciation'},
    @{Row=476; I='ba'; J='Appreciation'},
    @{Row=480; I='sd'; J='Statement-non-opinion'},
    @{Row=484; I='ba'; J='Appreciation'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
